$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 58233
$ws.Range("J3").Value = 58233
$ws.Range("L3").Value = 58233
$ws.Range("N3").Value = -58461

$ws.Range("H17").Value = 380.87878
$ws.Range("J17").Value = 319
$ws.Range("L17").Value = 957
$ws.Range("N17").Value = -1293

$ws.Range("H61").Value = 306.42856
$ws.Range("I61").Value = 292.5
$ws.Range("J61").Value = 390
$ws.Range("K61").Value = 877.5
$ws.Range("L61").Value = 1170
$ws.Range("M61").Value = -705.5
$ws.Range("N61").Value = -1514

$ws.Range("H70").Value = 1850
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 1850
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 5550
$ws.Range("M70").Value = ""
$ws.Range("N70").Value = -6090

$ws.Range("H73").Value = 1850
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 1850
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 5550
$ws.Range("M73").Value = ""
$ws.Range("N73").Value = -7422

$ws.Range("H94").Value = 957.73334
$ws.Range("I94").Value = 997.5714
$ws.Range("K94").Value = 997.5714
$ws.Range("M94").Value = -546.5714

$ws.Range("H102").Value = 58233
$ws.Range("J102").Value = 58233
$ws.Range("L102").Value = 58233
$ws.Range("N102").Value = -64723

$ws.Range("H115").Value = 744.6429000000001
$ws.Range("I115").Value = 744.6429000000001
$ws.Range("K115").Value = 2233.9287
$ws.Range("M115").Value = -666.9287000000004

$ws.Range("H132").Value = 10785.792
$ws.Range("I132").Value = 1065.1305
$ws.Range("J132").Value = 74664.42999999999
$ws.Range("K132").Value = 3195.3915
$ws.Range("L132").Value = 223993.29
$ws.Range("M132").Value = -665.3914999999997
$ws.Range("N132").Value = -229053.29

$ws.Range("H137").Value = 3741.2744
$ws.Range("I137").Value = 4117.5
$ws.Range("K137").Value = 12352.5
$ws.Range("M137").Value = -9802.5

$ws.Range("H141").Value = 9132.429
$ws.Range("I141").Value = 7077.636
$ws.Range("K141").Value = 21232.908
$ws.Range("M141").Value = -16052.908

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 31584
$ws.Range("I2").Value = 56835.168
$ws.Range("J2").Value = 6332.8335
$ws.Range("K2").Value = 56835.168
$ws.Range("L2").Value = 6332.8335
$ws.Range("M2").Value = -56722.168
$ws.Range("N2").Value = -6558.8335

$ws.Range("H61").Value = 2313.4546
$ws.Range("I61").Value = 1455.6471
$ws.Range("J61").Value = 5230
$ws.Range("K61").Value = 1455.6471
$ws.Range("L61").Value = 5230
$ws.Range("M61").Value = -1243.6471
$ws.Range("N61").Value = -5654

$ws.Range("H97").Value = 2867.923
$ws.Range("I97").Value = 1108
$ws.Range("K97").Value = 1108
$ws.Range("M97").Value = -612

$ws.Range("H110").Value = 5424.6665
$ws.Range("I110").Value = 5290.25
$ws.Range("J110").Value = 6500
$ws.Range("K110").Value = 5290.25
$ws.Range("L110").Value = 6500
$ws.Range("M110").Value = -3245.25
$ws.Range("N110").Value = -10590

$ws.Range("H116").Value = 31584
$ws.Range("I116").Value = 56835.168
$ws.Range("J116").Value = 6332.8335
$ws.Range("K116").Value = 56835.168
$ws.Range("L116").Value = 6332.8335
$ws.Range("M116").Value = -54541.168
$ws.Range("N116").Value = -10920.8335

$ws.Range("H124").Value = 79954.5
$ws.Range("J124").Value = 79954.5
$ws.Range("L124").Value = 79954.5
$ws.Range("N124").Value = -89774.5

$ws.Range("H132").Value = 1541.5652
$ws.Range("I132").Value = 1312.4286
$ws.Range("J132").Value = 3947.5
$ws.Range("K132").Value = 3937.2858
$ws.Range("L132").Value = 11842.5
$ws.Range("M132").Value = -1407.2858
$ws.Range("N132").Value = -16902.5

$ws.Range("H136").Value = 2313.4546
$ws.Range("I136").Value = 1455.6471
$ws.Range("J136").Value = 5230
$ws.Range("K136").Value = 4366.9413
$ws.Range("L136").Value = 15690
$ws.Range("M136").Value = -1816.9413
$ws.Range("N136").Value = -20790

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 31584
$ws.Range("I3").Value = 56835.168
$ws.Range("J3").Value = 6332.8335
$ws.Range("K3").Value = 56835.168
$ws.Range("L3").Value = 6332.8335
$ws.Range("M3").Value = -56721.168
$ws.Range("N3").Value = -6560.8335

$ws.Range("H20").Value = 6025
$ws.Range("J20").Value = 10000
$ws.Range("L20").Value = 10000
$ws.Range("N20").Value = -10494

$ws.Range("H94").Value = 1241.7693
$ws.Range("I94").Value = 417.57144
$ws.Range("J94").Value = 3339.7273
$ws.Range("K94").Value = 417.57144
$ws.Range("L94").Value = 3339.7273
$ws.Range("M94").Value = 33.42856
$ws.Range("N94").Value = -4241.7273

$ws.Range("H99").Value = 83095.60000000001
$ws.Range("I99").Value = 400010
$ws.Range("J99").Value = 3867
$ws.Range("K99").Value = 400010
$ws.Range("L99").Value = 3867
$ws.Range("M99").Value = -398512
$ws.Range("N99").Value = -6863

$ws.Range("H105").Value = 2867.074
$ws.Range("I105").Value = 2107.6667
$ws.Range("K105").Value = 2107.6667
$ws.Range("M105").Value = -360.6667000000002

$ws.Range("H107").Value = 14058.223
$ws.Range("I107").Value = 15251.833
$ws.Range("J107").Value = 11671
$ws.Range("K107").Value = 15251.833
$ws.Range("L107").Value = 11671
$ws.Range("M107").Value = -13331.833
$ws.Range("N107").Value = -15511

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 335
$ws.Range("I22").Value = 298.3
$ws.Range("J22").Value = 408.4
$ws.Range("K22").Value = 298.3
$ws.Range("L22").Value = 408.4
$ws.Range("M22").Value = 51.69999999999999
$ws.Range("N22").Value = -1108.4

$ws.Range("H31").Value = 1926.5862
$ws.Range("I31").Value = 1585.0454
$ws.Range("K31").Value = 1585.0454
$ws.Range("M31").Value = -1290.0454

$ws.Range("H34").Value = 1926.5862
$ws.Range("I34").Value = 1585.0454
$ws.Range("K34").Value = 1585.0454
$ws.Range("M34").Value = -1383.0454

$ws.Range("H93").Value = 19717.834
$ws.Range("I93").Value = 3661.4
$ws.Range("J93").Value = 100000
$ws.Range("K93").Value = 3661.4
$ws.Range("L93").Value = 100000
$ws.Range("M93").Value = -1789.4
$ws.Range("N93").Value = -103744

$ws.Range("H94").Value = 1944
$ws.Range("I94").Value = 1400
$ws.Range("K94").Value = 1400
$ws.Range("M94").Value = -949

$ws.Range("H99").Value = 2649.3572
$ws.Range("I99").Value = 2683.6924
$ws.Range("K99").Value = 2683.6924
$ws.Range("M99").Value = -1185.6924

$ws.Range("H126").Value = 2649.3572
$ws.Range("I126").Value = 2683.6924
$ws.Range("K126").Value = 8051.0772
$ws.Range("M126").Value = -5581.0772

$ws.Range("H132").Value = 3112.3809
$ws.Range("I132").Value = 2186.7222
$ws.Range("J132").Value = 8666.333000000001
$ws.Range("K132").Value = 6560.1666
$ws.Range("L132").Value = 25998.999
$ws.Range("M132").Value = -4030.1666
$ws.Range("N132").Value = -31058.999

$ws.Range("H134").Value = 1257
$ws.Range("I134").Value = 1163.8462
$ws.Range("K134").Value = 3491.5386
$ws.Range("M134").Value = -956.5385999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 55833948
$ws.Range("I4").Value = 59118004
$ws.Range("J4").Value = 5000
$ws.Range("K4").Value = 177354012
$ws.Range("L4").Value = 15000
$ws.Range("M4").Value = -177353900
$ws.Range("N4").Value = -15224

$ws.Range("H56").Value = 6498.75
$ws.Range("I56").Value = 6498.75
$ws.Range("K56").Value = 6498.75
$ws.Range("M56").Value = -5968.75

$ws.Range("H140").Value = 2980
$ws.Range("I140").Value = 2870
$ws.Range("K140").Value = 8610
$ws.Range("M140").Value = -3430

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 74443.31
$ws.Range("I80").Value = 280527.75
$ws.Range("J80").Value = 5748.5
$ws.Range("K80").Value = 280527.75
$ws.Range("L80").Value = 5748.5
$ws.Range("M80").Value = -279529.75
$ws.Range("N80").Value = -7744.5

$ws.Range("H83").Value = 74443.31
$ws.Range("I83").Value = 280527.75
$ws.Range("J83").Value = 5748.5
$ws.Range("K83").Value = 1402638.75
$ws.Range("L83").Value = 28742.5
$ws.Range("M83").Value = -1397646.75
$ws.Range("N83").Value = -38726.5

$ws.Range("H102").Value = 2205.4075
$ws.Range("I102").Value = 2001.7693
$ws.Range("J102").Value = 7500
$ws.Range("K102").Value = 2001.7693
$ws.Range("L102").Value = 7500
$ws.Range("M102").Value = -379.7692999999999
$ws.Range("N102").Value = -10744

$ws.Range("H107").Value = 374.83334
$ws.Range("I107").Value = 299.81818
$ws.Range("K107").Value = 299.81818
$ws.Range("M107").Value = 1620.18182

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1435.6
$ws.Range("I22").Value = 1063.3334
$ws.Range("J22").Value = 1994
$ws.Range("K22").Value = 1063.3334
$ws.Range("L22").Value = 1994
$ws.Range("M22").Value = -768.3334
$ws.Range("N22").Value = -2584

$ws.Range("H27").Value = 1435.6
$ws.Range("I27").Value = 1063.3334
$ws.Range("J27").Value = 1994
$ws.Range("K27").Value = 1063.3334
$ws.Range("L27").Value = 1994
$ws.Range("M27").Value = -956.3334
$ws.Range("N27").Value = -2208

$ws.Range("H46").Value = 2051.158
$ws.Range("I46").Value = 1188.75
$ws.Range("J46").Value = 2678.3635
$ws.Range("K46").Value = 1188.75
$ws.Range("L46").Value = 2678.3635
$ws.Range("M46").Value = -1000.75
$ws.Range("N46").Value = -3054.3635

$ws.Range("H122").Value = 5337
$ws.Range("I122").Value = 2736
$ws.Range("J122").Value = 6897.6
$ws.Range("K122").Value = 8208
$ws.Range("L122").Value = 20692.8
$ws.Range("M122").Value = -5758
$ws.Range("N122").Value = -25592.8

$ws.Range("H133").Value = 80000
$ws.Range("J133").Value = 80000
$ws.Range("L133").Value = 80000
$ws.Range("N133").Value = -85060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = ""
$ws.Range("N113").Value = ""

$ws.Range("H122").Value = 2689.5
$ws.Range("I122").Value = 2632.25
$ws.Range("K122").Value = 7896.75
$ws.Range("M122").Value = -5446.75

$ws.Range("H133").Value = 125000
$ws.Range("J133").Value = 125000
$ws.Range("L133").Value = 125000
$ws.Range("N133").Value = -135120
